{"js": "// Adiciona \"nome\" como requisito de cadastro do usu\u00e1rio (RF01).\n// Before: \", exigindo um e-mail e uma senha.\"\n// After : \", exigindo um e-mail (login), nome e senha.\"\nconst body = context.document.body;\n\nconst results = body.search(\", exigindo um e-mail e uma senha.\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target sentence not found: ', exigindo um e-mail e uma senha.'\");\n}\n\nresults.items[0].insertText(\n  \", exigindo um e-mail (login), nome e senha.\",\n  \"Replace\"\n);\nawait context.sync();\n", "ps1": "# Adiciona \"nome\" como requisito de cadastro do usu\u00e1rio (RF01).\n# Before: \", exigindo um e-mail e uma senha.\"\n# After : \", exigindo um e-mail (login), nome e senha.\"\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \", exigindo um e-mail e uma senha.\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \", exigindo um e-mail (login), nome e senha.\"\n\n# wdFindContinue = 1, wdReplaceAll = 2\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
